$d = $word.ActiveDocument

# Ordered list of (old, new) text replacements taken from the diff.
# A two-pass placeholder strategy is used so that a "new" value that happens
# to textually contain a not-yet-processed "old" value (or vice versa) cannot
# cause an accidental double replacement.
$pairs = @(
    ,@("2025-06-30 Monday", "2025-07-01 Tuesday")
    ,@("29+9=", "41+21=")
    ,@("41-39=", "12+87=")
    ,@("45-30=", "27-25=")
    ,@("2+29=", "39+16=")
    ,@("27+45=", "80-5=")
    ,@("89-80=", "55-31=")
    ,@("64+29=", "70-55=")
    ,@("62-48=", "59-38=")
    ,@("64-17=", "88-0=")
    ,@("5+72=", "77-17=")
    ,@("39+50=", "13+45=")
    ,@("91-28=", "76-59=")
    ,@("45+33=", "1+37=")
    ,@("98-59=", "34+55=")
    ,@("71+14=", "68-27=")
    ,@("6+8=", "13+55=")
    ,@("42+35=", "18+33=")
    ,@("39-38=", "31-17=")
    ,@("32+33=", "18+37=")
    ,@("9+26=", "73+25=")
    ,@("83-61=", "58-8=")
    ,@("56-44=", "73+2=")
    ,@("88-53=", "46+10=")
    ,@("0+42=", "4-3=")
    ,@("14+14=", "44-2=")
    ,@("56+14=", "68+5=")
    ,@("81-47=", "88-12=")
    ,@("10+44=", "67-46=")
    ,@("91+3=", "55+0=")
    ,@("25-13=", "66-65=")
    ,@("38+57=", "23+45=")
    ,@("48-11=", "83-31=")
    ,@("27-21=", "7+37=")
    ,@("23-17=", "94-48=")
    ,@("49-28=", "96+1=")
    ,@("88-22=", "55-46=")
    ,@("84-61=", "87-57=")
    ,@("68+1=", "45+23=")
    ,@("40+4=", "19+31=")
    ,@("57-11=", "90-57=")
    ,@("66-60=", "73-28=")
    ,@("72-27=", "62+15=")
    ,@("79+3=", "68+11=")
    ,@("47-2=", "83-1=")
    ,@("25-24=", "98-5=")
    ,@("65-22=", "45+22=")
    ,@("63-52=", "99-47=")
    ,@("67-21=", "6+38=")
    ,@("49+11=", "67-9=")
    ,@("63+26=", "97-67=")
    ,@("24+17=", "67+5=")
    ,@("56+10=", "78-31=")
    ,@("99-92=", "41+40=")
    ,@("23+3=", "63+17=")
    ,@("62+24=", "63-50=")
    ,@("41+49=", "29+52=")
    ,@("57-14=", "1+54=")
    ,@("94-35=", "72-48=")
    ,@("74-14=", "24+10=")
    ,@("15+1=", "85+8=")
    ,@("92-71=", "95-36=")
    ,@("67-64=", "74-34=")
    ,@("90-20=", "23+62=")
    ,@("11-2=", "26+71=")
    ,@("60+23=", "83-82=")
    ,@("82-66=", "43+41=")
    ,@("87-24=", "91-32=")
    ,@("76+4=", "99-27=")
    ,@("42-39=", "25-17=")
    ,@("25+18=", "50-6=")
    ,@("20+7=", "44-29=")
    ,@("95-16=", "95-54=")
    ,@("33-31=", "5+10=")
    ,@("96-43=", "32-29=")
    ,@("94-58=", "38-9=")
    ,@("88+11=", "96-54=")
    ,@("1+40=", "93-72=")
    ,@("99-15=", "17+81=")
    ,@("27+44=", "86-24=")
    ,@("5+18=", "54-47=")
    ,@("10+40=", "38-34=")
    ,@("31+22=", "2+60=")
    ,@("34-6=", "40+56=")
    ,@("40-26=", "37+38=")
    ,@("93-21=", "74-69=")
    ,@("65+28=", "21+65=")
    ,@("5+23=", "88-83=")
    ,@("15+31=", "27+50=")
    ,@("17+53=", "61-57=")
    ,@("63+25=", "91+7=")
    ,@("8+82=", "67-28=")
    ,@("60-53=", "78-32=")
    ,@("41-11=", "33-28=")
    ,@("81-31=", "54-2=")
    ,@("15+50=", "9-7=")
    ,@("10+59=", "0+95=")
    ,@("22+73=", "63-3=")
    ,@("17+64=", "72-9=")
    ,@("41+17=", "50+17=")
    ,@("59-10=", "7-5=")
)

# Pass 1: replace each original value with a unique, collision-free placeholder.
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $placeholder = "@@{0}@@" -f $i
    $d.Content.Find.Execute($pairs[$i][0], $false, $false, $false, $false, $false, $true, 1, $false, $placeholder, 2) | Out-Null
}

# Pass 2: replace each placeholder with the final new value.
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $placeholder = "@@{0}@@" -f $i
    $d.Content.Find.Execute($placeholder, $false, $false, $false, $false, $false, $true, 1, $false, $pairs[$i][1], 2) | Out-Null
}
